$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'GATES  9:11572,30060,GATES  9,MIDWAY10/11572,11579,GATES  9,GATES 16'
$ws.Range("A3").Value = 'GATES 10:11573,30050,GATES 10,LOSBNS 4/11573,11578,GATES 10,GATES 15'
$ws.Range("A4").Value = 'GATES 11:30056,11574,GATES  7,GATES 11/11574,11577,GATES 11,GATES 14'
$ws.Range("A5").Value = 'GATES 12:11575,11576,GATES 12,GATES 13/11575,30050,GATES 12,LOSBNS 4'
$ws.Range("A6").Value = 'GATES 13:30055,11576,GATES  8,GATES 13/11575,11576,GATES 12,GATES 13'
$ws.Range("A7").Value = 'GATES 14:11574,11577,GATES 11,GATES 14/11577,30060,GATES 14,MIDWAY10'
$ws.Range("A8").Value = 'GATES 15:30056,11578,GATES  7,GATES 15/11573,11578,GATES 10,GATES 15'
$ws.Range("A9").Value = 'GATES 16:30055,11579,GATES  8,GATES 16/11572,11579,GATES  9,GATES 16'
$ws.Range("A10").Value = 'GATES  4'
$ws.Range("A11").Value = 'METCLF 4'
$ws.Range("A12").Value = 'METCLF 5'
$ws.Range("A13").Value = 'MIDWAY11:30060,11747,MIDWAY10,MIDWAY11'
$ws.Range("A14").Value = 'MIDWAY12:30060,11748,MIDWAY10,MIDWAY12'
$ws.Range("A15").Value = 'MIDWAY13:30060,11749,MIDWAY10,MIDWAY13'
$ws.Range("A16").Value = 'MIDWAY 1'
$ws.Range("A17").Value = 'RNDMTN 5:11878,11949,RNDMTN 5,TBLMTN11/30005,11878,RNDMTN 4,RNDMTN 5'
$ws.Range("A18").Value = 'RNDMTN 6:11879,11882,RNDMTN 6,RNDMTN 9/30005,11879,RNDMTN 4,RNDMTN 6'
$ws.Range("A19").Value = 'RNDMTN 7:11880,11948,RNDMTN 7,TBLMTN10/30005,11880,RNDMTN 4,RNDMTN 7'
$ws.Range("A20").Value = 'RNDMTN 8:30005,11881,RNDMTN 4,RNDMTN 8/11881,11883,RNDMTN 8,RNDMTN10'
$ws.Range("A21").Value = 'RNDMTN 9:11879,11882,RNDMTN 6,RNDMTN 9'
$ws.Range("A22").Value = 'RNDMTN10:30010,11883,INDNSP 1,RNDMTN10/11881,11883,RNDMTN 8,RNDMTN10'
$ws.Range("A23").Value = 'TBLMTN 7:11945,11965,TBLMTN 7,TESLA 12/11945,11952,TBLMTN 7,TBLMTN12'
$ws.Range("A24").Value = 'TBLMTN 8:30015,11946,TBLMTN 6,TBLMTN 8/11946,11947,TBLMTN 8,TBLMTN 9'
$ws.Range("A25").Value = 'TBLMTN 9:11947,11990,TBLMTN 9,VACADX 9/11946,11947,TBLMTN 8,TBLMTN 9'
$ws.Range("A26").Value = 'TBLMTN10:11880,11948,RNDMTN 7,TBLMTN10/30015,11948,TBLMTN 6,TBLMTN10'
$ws.Range("A27").Value = 'TBLMTN11:11878,11949,RNDMTN 5,TBLMTN11/30015,11949,TBLMTN 6,TBLMTN11'
$ws.Range("A28").Value = 'TBLMTN 2'
$ws.Range("A29").Value = 'TBLMTN 3'
$ws.Range("A30").Value = 'TBLMTN12:30015,11952,TBLMTN 6,TBLMTN12/11945,11952,TBLMTN 7,TBLMTN12'
$ws.Range("A31").Value = 'TESLA 12:11945,11965,TBLMTN 7,TESLA 12/11965,11966,TESLA 12,TESLA 13'
$ws.Range("A32").Value = 'TESLA 13:30040,11966,TESLA 10,TESLA 13/11965,11966,TESLA 12,TESLA 13'
$ws.Range("A33").Value = 'VACADX 9:11947,11990,TBLMTN 9,VACADX 9/11990,11993,VACADX 9,VACADX12'
$ws.Range("A34").Value = 'VACADX10:30040,11991,TESLA 10,VACADX10/11991,11992,VACADX10,VACADX11'
$ws.Range("A35").Value = 'VACADX11:30030,11992,VACADX 8,VACADX11/11991,11992,VACADX10,VACADX11'
$ws.Range("A36").Value = 'VACADX12:11990,11993,VACADX 9,VACADX12/30030,11993,VACADX 8,VACADX12'
$ws.Range("A37").Value = 'RNDMTN 4:30005,11881,RNDMTN 4,RNDMTN 8/30005,11879,RNDMTN 4,RNDMTN 6/30005,11880,RNDMTN 4,RNDMTN 7/30005,11878,RNDMTN 4,RNDMTN 5'
$ws.Range("A38").Value = 'INDNSP 1:30010,11883,INDNSP 1,RNDMTN10'
$ws.Range("A39").Value = 'TBLMTN 6:30015,11952,TBLMTN 6,TBLMTN12/30015,11948,TBLMTN 6,TBLMTN10/30015,11946,TBLMTN 6,TBLMTN 8/30015,11949,TBLMTN 6,TBLMTN11'
$ws.Range("A40").Value = 'VACADX 8:30030,11992,VACADX 8,VACADX11/30030,11993,VACADX 8,VACADX12'
$ws.Range("A41").Value = 'TESLA 10:30040,11991,TESLA 10,VACADX10/30050,30040,LOSBNS 4,TESLA 10/30042,30040,METCLF 5,TESLA 10/30040,99006,TESLA 10,TESLA 11/30040,11966,TESLA 10,TESLA 13'
$ws.Range("A42").Value = 'METCLF 5:30042,30040,METCLF 5,TESLA 10/30042,30045,METCLF 5,MOSSLD13'
$ws.Range("A43").Value = 'MOSSLD13:30042,30045,METCLF 5,MOSSLD13/30045,30046,MOSSLD13,VISTRA 4/30050,30045,LOSBNS 4,MOSSLD13'
$ws.Range("A44").Value = 'VISTRA 4:30045,30046,MOSSLD13,VISTRA 4'
$ws.Range("A45").Value = 'LOSBNS 4:11573,30050,GATES 10,LOSBNS 4/11575,30050,GATES 12,LOSBNS 4/30050,30040,LOSBNS 4,TESLA 10/30055,30050,GATES  8,LOSBNS 4/30050,99005,LOSBNS 4,LOSBNS 5/30050,30045,LOSBNS 4,MOSSLD13'
$ws.Range("A46").Value = 'GATES  8:30057,30055,DIABLO 4,GATES  8/30055,30050,GATES  8,LOSBNS 4/30055,11576,GATES  8,GATES 13/30055,11579,GATES  8,GATES 16'
$ws.Range("A47").Value = 'GATES  7:30056,11578,GATES  7,GATES 15/30056,11574,GATES  7,GATES 11'
$ws.Range("A48").Value = 'DIABLO 4:30057,30060,DIABLO 4,MIDWAY10/30057,30055,DIABLO 4,GATES  8/30057,30060,DIABLO 4,MIDWAY10'
$ws.Range("A49").Value = 'MIDWAY10:30060,11749,MIDWAY10,MIDWAY13/11577,30060,GATES 14,MIDWAY10/30060,11748,MIDWAY10,MIDWAY12/11572,30060,GATES  9,MIDWAY10/30057,30060,DIABLO 4,MIDWAY10/30060,11747,MIDWAY10,MIDWAY11/30057,30060,DIABLO 4,MIDWAY10'
$ws.Range("A50").Value = 'LOSBNS 5:30050,99005,LOSBNS 4,LOSBNS 5'
$ws.Range("A51").Value = 'TESLA 11:30040,99006,TESLA 10,TESLA 11'
